$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the two competition result URLs (A7 and A18) to the new links.
$ws.Range("A7").Value = "https://www.rfebm.com/competiciones/resultados_completos.php?seleccion=0&id=1028231"
$ws.Range("A18").Value = "https://www.rfebm.com/competiciones/resultados_completos.php?seleccion=0&id=1029476"

# Match the styling already used on similar "plain text link" cells (A13/A15):
# underlined text without an actual hyperlink object.
$ws.Range("A7").Font.Underline = 1
$ws.Range("A18").Font.Underline = 1

# Move the active selection to A18, matching the saved cursor position.
$ws.Range("A18").Select()
